$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.673822
$ws.Range("H2").Value = 8.021466
$ws.Range("I2").Value = 0.8002273347603108
$ws.Range("J2").Value = 0.8002273347603108
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.847811666666666
$ws.Range("N2").Value = 11.543435
$ws.Range("O2").Value = 0.0396810199351781
$ws.Range("P2").Value = 0.03968101993517809
$ws.Range("Q2").Value = 10.28836348619
$ws.Range("R2").Value = 92.59527137571
$ws.Range("S2").Value = 0.03175383682329833
$ws.Range("T2").Value = 0.03175383682329833

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.673822
$ws.Range("H3").Value = 8.021466
$ws.Range("I3").Value = 0.8002273347603108
$ws.Range("J3").Value = 0.8002273347603108
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 57.80210733333333
$ws.Range("N3").Value = 173.406322
$ws.Range("O3").Value = 0.5960911739155557
$ws.Range("P3").Value = 0.5960911739155557
$ws.Range("Q3").Value = 154.552546234228
$ws.Range("R3").Value = 1390.972916108052
$ws.Range("S3").Value = 0.47700845137659
$ws.Range("T3").Value = 0.47700845137659

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.673822
$ws.Range("H4").Value = 8.021466
$ws.Range("I4").Value = 0.8002273347603108
$ws.Range("J4").Value = 0.8002273347603108
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.066157
$ws.Range("N4").Value = 18.198471
$ws.Range("O4").Value = 0.06255797260873913
$ws.Range("P4").Value = 0.06255797260873913
$ws.Range("Q4").Value = 16.219824042054
$ws.Range("R4").Value = 145.978416378486
$ws.Range("S4").Value = 0.05006059968869984
$ws.Range("T4").Value = 0.05006059968869984

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.673822
$ws.Range("H5").Value = 8.021466
$ws.Range("I5").Value = 0.8002273347603108
$ws.Range("J5").Value = 0.8002273347603108
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.25249166666667
$ws.Range("N5").Value = 87.757475
$ws.Range("O5").Value = 0.3016698335405271
$ws.Range("P5").Value = 0.301669833540527
$ws.Range("Q5").Value = 78.21595577315
$ws.Range("R5").Value = 703.94360195835
$ws.Range("S5").Value = 0.2414044468717226
$ws.Range("T5").Value = 0.2414044468717226

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.667506
$ws.Range("H6").Value = 2.002518
$ws.Range("I6").Value = 0.1997726652396891
$ws.Range("J6").Value = 0.1997726652396891
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.847811666666666
$ws.Range("N6").Value = 11.543435
$ws.Range("O6").Value = 0.0396810199351781
$ws.Range("P6").Value = 0.03968101993517809
$ws.Range("Q6").Value = 2.56843737437
$ws.Range("R6").Value = 23.11593636933
$ws.Range("S6").Value = 0.007927183111879767
$ws.Range("T6").Value = 0.007927183111879765

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.667506
$ws.Range("H7").Value = 2.002518
$ws.Range("I7").Value = 0.1997726652396891
$ws.Range("J7").Value = 0.1997726652396891
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 57.80210733333333
$ws.Range("N7").Value = 173.406322
$ws.Range("O7").Value = 0.5960911739155557
$ws.Range("P7").Value = 0.5960911739155557
$ws.Range("Q7").Value = 38.583253457644
$ws.Range("R7").Value = 347.249281118796
$ws.Range("S7").Value = 0.1190827225389656
$ws.Range("T7").Value = 0.1190827225389656

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.667506
$ws.Range("H8").Value = 2.002518
$ws.Range("I8").Value = 0.1997726652396891
$ws.Range("J8").Value = 0.1997726652396891
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.066157
$ws.Range("N8").Value = 18.198471
$ws.Range("O8").Value = 0.06255797260873913
$ws.Range("P8").Value = 0.06255797260873913
$ws.Range("Q8").Value = 4.049196194442001
$ws.Range("R8").Value = 36.44276574997801
$ws.Range("S8").Value = 0.01249737292003929
$ws.Range("T8").Value = 0.01249737292003929

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.667506
$ws.Range("H9").Value = 2.002518
$ws.Range("I9").Value = 0.1997726652396891
$ws.Range("J9").Value = 0.1997726652396891
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 29.25249166666667
$ws.Range("N9").Value = 87.757475
$ws.Range("O9").Value = 0.3016698335405271
$ws.Range("P9").Value = 0.301669833540527
$ws.Range("Q9").Value = 19.52621370245
$ws.Range("R9").Value = 175.73592332205
$ws.Range("S9").Value = 0.06026538666880446
$ws.Range("T9").Value = 0.06026538666880445
